$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 194.42857
$ws.Range("I11").Value = 194.42857
$ws.Range("K11").Value = 194.42857
$ws.Range("M11").Value = -54.42857000000001

$ws.Range("H32").Value = 3674.1667
$ws.Range("I32").Value = 2874.75
$ws.Range("J32").Value = 4073.875
$ws.Range("K32").Value = 2874.75
$ws.Range("L32").Value = 4073.875
$ws.Range("M32").Value = -2548.75
$ws.Range("N32").Value = -4725.875

$ws.Range("H39").Value = 2601.7334
$ws.Range("I39").Value = 704.625
$ws.Range("J39").Value = 4769.857
$ws.Range("K39").Value = 2113.875
$ws.Range("L39").Value = 14309.571
$ws.Range("M39").Value = -1817.875
$ws.Range("N39").Value = -14901.571

$ws.Range("H43").Value = 7324.75
$ws.Range("J43").Value = 8933
$ws.Range("L43").Value = 8933
$ws.Range("N43").Value = -9071

$ws.Range("H121").Value = 1061.3529
$ws.Range("J121").Value = 1061.3529
$ws.Range("L121").Value = 3184.0587
$ws.Range("N121").Value = -6678.0587

$ws.Range("H135").Value = 3370.0833
$ws.Range("I135").Value = 3446.7368
$ws.Range("K135").Value = 31020.6312
$ws.Range("M135").Value = -28485.6312

$ws.Range("H138").Value = 10551.301
$ws.Range("J138").Value = 10730.831
$ws.Range("L138").Value = 32192.493
$ws.Range("N138").Value = -42472.493

$ws.Range("H141").Value = 2870.9092
$ws.Range("I141").Value = 2928.5
$ws.Range("K141").Value = 8785.5
$ws.Range("M141").Value = -3605.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 50000
$ws.Range("J22").Value = 50000
$ws.Range("L22").Value = 50000
$ws.Range("N22").Value = -50598

$ws.Range("H32").Value = 9037.875
$ws.Range("I32").Value = 6640.4
$ws.Range("K32").Value = 6640.4
$ws.Range("M32").Value = -6353.4

$ws.Range("H102").Value = 2529.5386
$ws.Range("I102").Value = 2407
$ws.Range("K102").Value = 2407
$ws.Range("M102").Value = -785

$ws.Range("H133").Value = 91998.8
$ws.Range("J133").Value = 91998.8
$ws.Range("L133").Value = 91998.8
$ws.Range("N133").Value = -97058.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1089.5
$ws.Range("I94").Value = 1205.25
$ws.Range("J94").Value = 742.25
$ws.Range("K94").Value = 1205.25
$ws.Range("L94").Value = 742.25
$ws.Range("M94").Value = -754.25
$ws.Range("N94").Value = -1644.25

$ws.Range("H122").Value = 89166.5
$ws.Range("J122").Value = 89166.5
$ws.Range("L122").Value = 89166.5
$ws.Range("N122").Value = -98966.5

$ws.Range("H126").Value = 94038.38
$ws.Range("J126").Value = 94038.38
$ws.Range("L126").Value = 94038.38
$ws.Range("N126").Value = -103918.38

$ws.Range("H130").Value = 88672.71000000001

$ws.Range("H132").Value = 89837.61
$ws.Range("J132").Value = 89837.61
$ws.Range("L132").Value = 89837.61
$ws.Range("N132").Value = -99957.61

$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

$ws.Range("H137").Value = 79701.53999999999
$ws.Range("J137").Value = 79701.53999999999
$ws.Range("L137").Value = 79701.53999999999
$ws.Range("N137").Value = -89901.53999999999

$ws.Range("H138").Value = 88186.30499999999
$ws.Range("J138").Value = 88186.30499999999
$ws.Range("L138").Value = 88186.30499999999
$ws.Range("N138").Value = -98466.30499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H58").Value = 6335.4443
$ws.Range("I58").Value = 6888
$ws.Range("J58").Value = 4898.8
$ws.Range("K58").Value = 6888
$ws.Range("L58").Value = 4898.8
$ws.Range("M58").Value = -6685
$ws.Range("N58").Value = -5304.8

$ws.Range("H86").Value = 25656026
$ws.Range("I86").Value = 33351444
$ws.Range("K86").Value = 33351444
$ws.Range("M86").Value = -33350321

$ws.Range("H89").Value = 25656026
$ws.Range("I89").Value = 33351444
$ws.Range("K89").Value = 166757220
$ws.Range("M89").Value = -166751604

$ws.Range("H136").Value = 6335.4443
$ws.Range("I136").Value = 6888
$ws.Range("J136").Value = 4898.8
$ws.Range("K136").Value = 20664
$ws.Range("L136").Value = 14696.4
$ws.Range("M136").Value = -18114
$ws.Range("N136").Value = -19796.4

$ws.Range("H141").Value = 287543.9
$ws.Range("J141").Value = 300046.22
$ws.Range("L141").Value = 300046.22
$ws.Range("N141").Value = -310406.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 6166.3335
$ws.Range("J9").Value = 8499.5
$ws.Range("L9").Value = 25498.5
$ws.Range("N9").Value = -25946.5

$ws.Range("H14").Value = 1162.8077
$ws.Range("I14").Value = 1162.8077
$ws.Range("K14").Value = 3488.4231
$ws.Range("M14").Value = -3315.4231

$ws.Range("H33").Value = 122.8
$ws.Range("J33").Value = 185.85715
$ws.Range("L33").Value = 1115.1429
$ws.Range("N33").Value = -1681.1429

$ws.Range("H63").Value = 396
$ws.Range("I63").Value = 396
$ws.Range("K63").Value = 1188
$ws.Range("M63").Value = -439

$ws.Range("H64").Value = 3993
$ws.Range("I64").Value = 1657.6666
$ws.Range("K64").Value = 4972.9998
$ws.Range("M64").Value = -4702.9998

$ws.Range("H66").Value = 396
$ws.Range("I66").Value = 396
$ws.Range("K66").Value = 3564
$ws.Range("M66").Value = 180

$ws.Range("H67").Value = 3993
$ws.Range("I67").Value = 1657.6666
$ws.Range("K67").Value = 4972.9998
$ws.Range("M67").Value = -4036.9998

$ws.Range("H118").Value = 5496.6665
$ws.Range("I118").Value = 4796
$ws.Range("K118").Value = 14388
$ws.Range("M118").Value = -13145

$ws.Range("H121").Value = 251434.75
$ws.Range("I121").Value = 1869.5
$ws.Range("J121").Value = 501000
$ws.Range("K121").Value = 5608.5
$ws.Range("L121").Value = 1503000
$ws.Range("M121").Value = -4298.5
$ws.Range("N121").Value = -1505620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7556.731
$ws.Range("I80").Value = 5879.4
$ws.Range("J80").Value = 8605.0625
$ws.Range("K80").Value = 5879.4
$ws.Range("L80").Value = 8605.0625
$ws.Range("M80").Value = -4881.4
$ws.Range("N80").Value = -10601.0625

$ws.Range("H83").Value = 7556.731
$ws.Range("I83").Value = 5879.4
$ws.Range("J83").Value = 8605.0625
$ws.Range("K83").Value = 29397
$ws.Range("L83").Value = 43025.3125
$ws.Range("M83").Value = -24405
$ws.Range("N83").Value = -53009.3125

$ws.Range("H107").Value = 386.75
$ws.Range("J107").Value = 1049
$ws.Range("L107").Value = 1049
$ws.Range("N107").Value = -4889

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13793
$ws.Range("J7").Value = 5217.143
$ws.Range("L7").Value = 5217.143
$ws.Range("N7").Value = -5441.143

$ws.Range("H61").Value = 1889.2778
$ws.Range("I61").Value = 2294.8333
$ws.Range("J61").Value = 1078.1666
$ws.Range("K61").Value = 2294.8333
$ws.Range("L61").Value = 1078.1666
$ws.Range("M61").Value = -2092.8333
$ws.Range("N61").Value = -1482.1666

$ws.Range("H113").Value = 1889.2778
$ws.Range("I113").Value = 2294.8333
$ws.Range("J113").Value = 1078.1666
$ws.Range("K113").Value = 2294.8333
$ws.Range("L113").Value = 1078.1666
$ws.Range("M113").Value = -124.8332999999998
$ws.Range("N113").Value = -5418.1666

$ws.Range("H126").Value = 13793
$ws.Range("J126").Value = 5217.143
$ws.Range("L126").Value = 15651.429
$ws.Range("N126").Value = -20591.429

$ws.Range("H132").Value = 16592.375
$ws.Range("I132").Value = 16737.066
$ws.Range("K132").Value = 50211.198
$ws.Range("M132").Value = -47681.198

$ws.Range("H136").Value = 8656.857
$ws.Range("I136").Value = 7287.4165
$ws.Range("J136").Value = 10482.777
$ws.Range("K136").Value = 21862.2495
$ws.Range("L136").Value = 31448.331
$ws.Range("M136").Value = -19312.2495
$ws.Range("N136").Value = -36548.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 71436.55499999999
$ws.Range("I75").Value = 35799.8
$ws.Range("J75").Value = 115982.5
$ws.Range("K75").Value = 35799.8
$ws.Range("L75").Value = 115982.5
$ws.Range("M75").Value = -34863.8
$ws.Range("N75").Value = -117854.5

$ws.Range("H78").Value = 71436.55499999999
$ws.Range("I78").Value = 35799.8
$ws.Range("J78").Value = 115982.5
$ws.Range("K78").Value = 107399.4
$ws.Range("L78").Value = 347947.5
$ws.Range("M78").Value = -102719.4
$ws.Range("N78").Value = -357307.5

$ws.Range("H118").Value = 74800
$ws.Range("I118").Value = 74800
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 74800
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -73143
$ws.Range("N118").ClearContents()
